$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 4099.909
$ws.Cells.Item(64, 9).Value = 3928.5
$ws.Cells.Item(64, 10).Value = 4399.875
$ws.Cells.Item(64, 11).Value = 3928.5
$ws.Cells.Item(64, 12).Value = 4399.875
$ws.Cells.Item(64, 13).Value = -3680.5
$ws.Cells.Item(64, 14).Value = -4895.875

$ws.Cells.Item(67, 8).Value = 4099.909
$ws.Cells.Item(67, 9).Value = 3928.5
$ws.Cells.Item(67, 10).Value = 4399.875
$ws.Cells.Item(67, 11).Value = 3928.5
$ws.Cells.Item(67, 12).Value = 4399.875
$ws.Cells.Item(67, 13).Value = -3070.5
$ws.Cells.Item(67, 14).Value = -6115.875

$ws.Cells.Item(76, 8).Value = 3181.818
$ws.Cells.Item(76, 9).Value = 3180
$ws.Cells.Item(76, 10).Value = 3200
$ws.Cells.Item(76, 11).Value = 3180
$ws.Cells.Item(76, 12).Value = 3200
$ws.Cells.Item(76, 13).Value = -2865
$ws.Cells.Item(76, 14).Value = -3830

$ws.Cells.Item(79, 8).Value = 3181.818
$ws.Cells.Item(79, 9).Value = 3180
$ws.Cells.Item(79, 10).Value = 3200
$ws.Cells.Item(79, 11).Value = 3180
$ws.Cells.Item(79, 12).Value = 3200
$ws.Cells.Item(79, 13).Value = -2088
$ws.Cells.Item(79, 14).Value = -5384

$ws.Cells.Item(92, 8).Value = 998.6
$ws.Cells.Item(92, 9).Value = 983.64703
$ws.Cells.Item(92, 11).Value = 983.64703
$ws.Cells.Item(92, 13).Value = 264.35297

$ws.Cells.Item(129, 8).Value = 849.53845
$ws.Cells.Item(129, 10).Value = 849.52
$ws.Cells.Item(129, 12).Value = 2548.56
$ws.Cells.Item(129, 14).Value = -12548.56

$ws.Cells.Item(138, 8).Value = 1668
$ws.Cells.Item(138, 10).Value = 4220.75
$ws.Cells.Item(138, 12).Value = 12662.25
$ws.Cells.Item(138, 14).Value = -22942.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2221.7048
$ws.Cells.Item(32, 9).Value = 1628.6604
$ws.Cells.Item(32, 10).Value = 6150.625
$ws.Cells.Item(32, 11).Value = 1628.6604
$ws.Cells.Item(32, 12).Value = 6150.625
$ws.Cells.Item(32, 13).Value = -1341.6604
$ws.Cells.Item(32, 14).Value = -6724.625

$ws.Cells.Item(45, 8).Value = 3466.4062
$ws.Cells.Item(45, 10).Value = 3209.682
$ws.Cells.Item(45, 12).Value = 3209.682
$ws.Cells.Item(45, 14).Value = -3963.682

$ws.Cells.Item(80, 8).Value = 47303
$ws.Cells.Item(80, 10).Value = 47303
$ws.Cells.Item(80, 12).Value = 47303
$ws.Cells.Item(80, 14).Value = -49299

$ws.Cells.Item(83, 8).Value = 47303
$ws.Cells.Item(83, 10).Value = 47303
$ws.Cells.Item(83, 12).Value = 141909
$ws.Cells.Item(83, 14).Value = -151893

$ws.Cells.Item(102, 8).Value = 5379.8335
$ws.Cells.Item(102, 9).Value = 3000
$ws.Cells.Item(102, 10).Value = 5855.8
$ws.Cells.Item(102, 11).Value = 3000
$ws.Cells.Item(102, 12).Value = 5855.8
$ws.Cells.Item(102, 13).Value = -1378
$ws.Cells.Item(102, 14).Value = -9099.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 2042.25
$ws.Cells.Item(105, 9).Value = 1998.6666
$ws.Cells.Item(105, 11).Value = 1998.6666
$ws.Cells.Item(105, 13).Value = -251.6666

$ws.Cells.Item(134, 8).Value = 4467.278
$ws.Cells.Item(134, 9).Value = 4553.5884
$ws.Cells.Item(134, 11).Value = 13660.7652
$ws.Cells.Item(134, 13).Value = -11125.7652

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 5600.364
$ws.Cells.Item(62, 9).Value = 3684
$ws.Cells.Item(62, 10).Value = 7900
$ws.Cells.Item(62, 11).Value = 3684
$ws.Cells.Item(62, 12).Value = 7900
$ws.Cells.Item(62, 13).Value = -3060
$ws.Cells.Item(62, 14).Value = -9148

$ws.Cells.Item(65, 8).Value = 5600.364
$ws.Cells.Item(65, 9).Value = 3684
$ws.Cells.Item(65, 10).Value = 7900
$ws.Cells.Item(65, 11).Value = 18420
$ws.Cells.Item(65, 12).Value = 39500
$ws.Cells.Item(65, 13).Value = -15300
$ws.Cells.Item(65, 14).Value = -45740

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(80, 8).Value = 3000
$ws.Cells.Item(80, 10).Value = 3000
$ws.Cells.Item(80, 12).Value = 9000
$ws.Cells.Item(80, 14).Value = -10872

$ws.Cells.Item(83, 8).Value = 3000
$ws.Cells.Item(83, 10).Value = 3000
$ws.Cells.Item(83, 12).Value = 27000
$ws.Cells.Item(83, 14).Value = -36360

$ws.Cells.Item(131, 8).Value = 724.97
$ws.Cells.Item(131, 10).Value = 747.23157
$ws.Cells.Item(131, 12).Value = 2241.69471
$ws.Cells.Item(131, 14).Value = -12321.69471

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4622.8
$ws.Cells.Item(70, 9).Value = 4408.2856
$ws.Cells.Item(70, 10).Value = 4738.3076
$ws.Cells.Item(70, 11).Value = 4408.2856
$ws.Cells.Item(70, 12).Value = 4738.3076
$ws.Cells.Item(70, 13).Value = -4138.2856
$ws.Cells.Item(70, 14).Value = -5278.3076

$ws.Cells.Item(73, 8).Value = 4622.8
$ws.Cells.Item(73, 9).Value = 4408.2856
$ws.Cells.Item(73, 10).Value = 4738.3076
$ws.Cells.Item(73, 11).Value = 4408.2856
$ws.Cells.Item(73, 12).Value = 4738.3076
$ws.Cells.Item(73, 13).Value = -3472.2856
$ws.Cells.Item(73, 14).Value = -6610.3076

$ws.Cells.Item(80, 8).Value = 3807.647
$ws.Cells.Item(80, 9).Value = 3405
$ws.Cells.Item(80, 10).Value = 4382.857
$ws.Cells.Item(80, 11).Value = 3405
$ws.Cells.Item(80, 12).Value = 4382.857
$ws.Cells.Item(80, 13).Value = -2407
$ws.Cells.Item(80, 14).Value = -6378.857

$ws.Cells.Item(83, 8).Value = 3807.647
$ws.Cells.Item(83, 9).Value = 3405
$ws.Cells.Item(83, 10).Value = 4382.857
$ws.Cells.Item(83, 11).Value = 17025
$ws.Cells.Item(83, 12).Value = 21914.285
$ws.Cells.Item(83, 13).Value = -12033
$ws.Cells.Item(83, 14).Value = -31898.285

$ws.Cells.Item(96, 8).Value = 19474
$ws.Cells.Item(96, 10).Value = 19474
$ws.Cells.Item(96, 12).Value = 19474
$ws.Cells.Item(96, 14).Value = -24966

$ws.Cells.Item(132, 8).Value = 29092.95
$ws.Cells.Item(132, 9).Value = 4569.1816
$ws.Cells.Item(132, 11).Value = 13707.5448
$ws.Cells.Item(132, 13).Value = -11177.5448

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 3674.125
$ws.Cells.Item(93, 9).Value = 3770.4285
$ws.Cells.Item(93, 10).Value = 3000
$ws.Cells.Item(93, 11).Value = 3770.4285
$ws.Cells.Item(93, 12).Value = 3000
$ws.Cells.Item(93, 13).Value = -2522.4285
$ws.Cells.Item(93, 14).Value = -5496

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(95, 8).Value = 32000
$ws.Cells.Item(95, 10).Value = 32000
$ws.Cells.Item(95, 12).Value = 32000
$ws.Cells.Item(95, 14).Value = -37492

$ws.Cells.Item(107, 8).Value = 3247905
$ws.Cells.Item(107, 9).Value = 328.125
$ws.Cells.Item(107, 11).Value = 984.375
$ws.Cells.Item(107, 13).Value = 935.625

$ws.Cells.Item(126, 8).Value = 1342.875
$ws.Cells.Item(126, 9).Value = 845.625
$ws.Cells.Item(126, 10).Value = 1840.125
$ws.Cells.Item(126, 11).Value = 2536.875
$ws.Cells.Item(126, 12).Value = 5520.375
$ws.Cells.Item(126, 13).Value = -66.875
$ws.Cells.Item(126, 14).Value = -10460.375

$ws.Cells.Item(132, 8).Value = 2630.7
$ws.Cells.Item(132, 9).Value = 2008.4286
$ws.Cells.Item(132, 11).Value = 6025.2858
$ws.Cells.Item(132, 13).Value = -3495.2858
